# The workbook under test has a "Test Suite" sheet where column C holds
# Y/N run flags (as shared strings). This commit flips the remaining "N"
# rows (C3:C7) to "Y" so that all test cases run, which also makes the
# "N" shared string unused (it will be dropped from sharedStrings.xml on
# save). The active selection also moves to C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C3:C7").Value = "Y"

$ws.Activate()
$ws.Range("C6").Select()
